$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1. Rename header row (row 1): Spanish labels -> snake_case codes
# -----------------------------------------------------------------
$ws.Range("A1").Value2 = "mx_state"
$ws.Range("B1").Value2 = "mx_municipality"
$ws.Range("C1").Value2 = "n_matriculas"
$ws.Range("D1").Value2 = "pct_matriculas"

# -----------------------------------------------------------------
# 2. Title-case the Spanish state/municipality names in columns A/B
#    ("de" -> "De", "del" -> "Del", "la" -> "La", "el" -> "El",
#     "los" -> "Los", "y" -> "Y", ...), matching .title()-style
#    capitalization, via .NET's invariant TextInfo.ToTitleCase.
# -----------------------------------------------------------------
$ti = [System.Globalization.CultureInfo]::InvariantCulture.TextInfo

function TitleCaseEs($s) {
    return $ti.ToTitleCase($s.ToLower())
}

for ($r = 2; $r -le 376; $r++) {
    foreach ($col in @("A","B")) {
        # A92 is handled explicitly below (it has an embedded CR/LF
        # in the source data that needs scrubbing, not just
        # title-casing, and re-writing it through the generic
        # title-case path would leave the row's auto-fit height
        # artifact behind even after the text is fixed).
        if ($r -eq 92 -and $col -eq "A") { continue }
        $cell = $ws.Range("$col$r")
        $v = $cell.Value2
        if ($v -ne $null) {
            $cell.Value2 = TitleCaseEs($v)
        }
    }
}

# -----------------------------------------------------------------
# 3. Row 92 col A had a stray CR + trailing newline baked into the
#    inline string ("Estado de México_x000D_\n"). Clean it up to the
#    plain, title-cased state name.
# -----------------------------------------------------------------
$ws.Range("A92").Value2 = "Estado De México"

# -----------------------------------------------------------------
# 4. A handful of D-column percentages (all the rows whose numerator
#    was 1, i.e. 1/1101) get recalculated to one ULP higher due to
#    the refreshed floating point division.
# -----------------------------------------------------------------
$target = 0.0009082652134423251
$replacement = 0.0009082652134423252

for ($r = 2; $r -le 376; $r++) {
    $cell = $ws.Range("D$r")
    $v = $cell.Value2
    if ($v -ne $null -and $v -eq $target) {
        $cell.Value2 = $replacement
    }
}

# -----------------------------------------------------------------
# 5. Drop the trailing footnote/source rows (378-382); the sheet's
#    real data ends at row 376 ("Total").
# -----------------------------------------------------------------
$ws.Rows("378:382").Delete() | Out-Null
